# Line Item Template.xlsx - 23 Nov 23 Development
# - Remove the "Tax rate" column (column AC) from the header row/table.
# - Remove the leftover formatted-but-empty rows (4:42) below the header.
# - Re-point the header row's cell style onto the (now-first) duplicate
#   wrap-text style so the redundant style entry is no longer referenced.
# - Reset the sheet view (no frozen/scrolled "topLeftCell", selection at A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray formatted rows below the header (rows 4 through 42) that
# only carried a wrap-text style on column C with no real data.
$ws.Rows("4:42").Delete()

# Delete the "Tax rate" column (AC) entirely; everything to its right
# (Non-deductible tax rate ... Coding Block) shifts one column left.
$ws.Columns("AC:AC").Delete()

# Re-apply wrap text to the header row so it binds to the existing
# wrap-text style instead of keeping its own separate (duplicate) style.
$ws.Range("A1:AL1").WrapText = $true

# Reset selection/view to A2, clearing the old scrolled-to-AB1 / AK1 state.
$ws.Range("A2").Select() | Out-Null
